# Slide 9 ("Wyniki badań heurystyk A*") - split the title into three runs:
#   "Wyniki badań " + "- heurystyki " + "A*"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Title
$tr = $sh.TextFrame.TextRange

$tr.Text = "Wyniki badań "
[void]$tr.InsertAfter("- heurystyki ")
[void]$tr.InsertAfter("A*")
